$d = $word.ActiveDocument

# --- Edit 1: split the final run of the Employee.java paragraph (para 25) ---
# "... Visible changes of code in EBeforePrivate & EAfterPrivate."
# splits " & EAfterPrivate." into " & " + EAfterPrivate(spell-wrapped) + "."
$employeeXml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">I made all the methods inside the Employee Class with no modifier. The methods are used externally so if I made them private there would be errors when handling the Employee. Making them as </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>they're called 'package-private' means they can only be accessed within this package which is exactly what we need. There is still a level of encapsulation.</w:t></w:r><w:r><w:t xml:space="preserve"> Visible changes of code in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EBeforePrivate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EAfterPrivate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(25).Range.InsertXML($employeeXml)

# --- Edit 2: insert the new documentation paragraphs after paragraph 29 ---
# (the first empty <w:p/> following the RandomAccessEmployeeRecord.java section)
$anchor = $d.Paragraphs.Item(29)
for ($i = 0; $i -lt 12; $i++) {
    $anchor.Range.InsertParagraphAfter()
}

$p30 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>RandomFile.java</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(30).Range.InsertXML($p30)
$p31 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">There is duplicate code in this file for closing a  file whether it's read or write. I renamed the method to be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>closeFile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() and everywhere where there was the duplicate code I called that m</w:t></w:r><w:r><w:t xml:space="preserve">ethod. I had to change the method calls in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EmployeeDetails</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. I showed this in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RFBeforeCloseFile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RFAfterClose</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(31).Range.InsertXML($p31)
$p32 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">I also changed the methods to be private if they're used inside the class and no modifier if </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>theyre</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> used outside the class inside the package.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(32).Range.InsertXML($p32)
$p33 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">In the method </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>changeRecords</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() he was declaring a new variable </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>currentRecordStart</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oldDetails</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> equal to the parameters which is unnecessary</w:t></w:r><w:r><w:t xml:space="preserve"> and bad practise</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> It is quicker and more direct to use the variable name given as a parameter instead. You can use these parameters in the functions and the code is now cleaner.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(33).Range.InsertXML($p33)
$p34 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Some methods need to be changed to have no modifier and he has empty catch blocks so I just added a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>system.out.println</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to display the error.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(34).Range.InsertXML($p34)
$p35 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>I have some of the corrections I</w:t></w:r><w:r><w:t xml:space="preserve"> made in the pictures </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RFOne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RFOneA</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RFOne</w:t></w:r><w:r><w:t>B</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(35).Range.InsertXML($p35)
# paragraph 36 stays empty
$p37 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>SearchByIdDialog.java</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(37).Range.InsertXML($p37)
$p38 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Changed the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>classname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and methods to have no modifier. Changed variables to be of private as they are only used inside the class. Screenshots called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SIDBefore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SIDAfter</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(38).Range.InsertXML($p38)
# paragraph 39 stays empty
$p40 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>SearchBySurnameDialog.java</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(40).Range.InsertXML($p40)
$p41 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Changed the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>classname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and methods to have no modifier. Changed variables to be of private as they are only used inside the class. Screenshots called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SBSBefore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SBSAfter</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs.Item(41).Range.InsertXML($p41)

Write-Output "final paragraph count: $($d.Paragraphs.Count)"
